$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 16.70093436015432
$ws.Range("C2").Value = 6.468860745156574
$ws.Range("D2").Value = 11.35442112295005
$ws.Range("F2").Value = 56.9220229332803
$ws.Range("G2").Value = 3.769390396562137
$ws.Range("K2").Value = 12.82697287926648
$ws.Range("L2").Value = 10.87050959397311
$ws.Range("M2").Value = 16.80557504187401
$ws.Range("B3").Value = 16.67055818930364
$ws.Range("C3").Value = 6.389049573503332
$ws.Range("D3").Value = 11.21209189896101
$ws.Range("F3").Value = 55.79519193304199
$ws.Range("G3").Value = 3.772983881830482
$ws.Range("K3").Value = 12.83962088372387
$ws.Range("L3").Value = 10.87471791188447
$ws.Range("M3").Value = 16.82996393019905
$ws.Range("B4").Value = 16.65813284678714
$ws.Range("C4").Value = 6.33780158523264
$ws.Range("D4").Value = 11.1226030171463
$ws.Range("F4").Value = 55.09314430096277
$ws.Range("G4").Value = 3.775302029571474
$ws.Range("K4").Value = 12.8525638025309
$ws.Range("L4").Value = 10.87881984688133
$ws.Range("M4").Value = 16.84859476777214
$ws.Range("B5").Value = 16.65463849258699
$ws.Range("C5").Value = 6.316352401624684
$ws.Range("D5").Value = 11.0856220458491
$ws.Range("F5").Value = 54.80476950130085
$ws.Range("G5").Value = 3.776274900733167
$ws.Range("K5").Value = 12.85913617042277
$ws.Range("L5").Value = 10.88087337882188
$ws.Range("M5").Value = 16.85710536122476
$ws.Range("B6").Value = 16.65415309247614
$ws.Range("C6").Value = 6.31275661408657
$ws.Range("D6").Value = 11.07945078449055
$ws.Range("F6").Value = 54.75675522482839
$ws.Range("G6").Value = 3.776438152144063
$ws.Range("K6").Value = 12.86030577417571
$ws.Range("L6").Value = 10.88123744098689
$ws.Range("M6").Value = 16.85857397386273
$ws.Range("B7").Value = 16.65807936437523
$ws.Range("C7").Value = 6.337514603842494
$ws.Range("D7").Value = 11.12210633922032
$ws.Range("F7").Value = 55.08926406213468
$ws.Range("G7").Value = 3.775315035697719
$ws.Range("K7").Value = 12.85264718995565
$ws.Range("L7").Value = 10.8788459947364
$ws.Range("M7").Value = 16.84870582761633
$ws.Range("B8").Value = 16.6891712130853
$ws.Range("C8").Value = 6.441806863606857
$ws.Range("D8").Value = 11.30578982544471
$ws.Range("F8").Value = 56.53576269162468
$ws.Range("G8").Value = 3.770606306970365
$ws.Range("K8").Value = 12.83025746707398
$ws.Range("L8").Value = 10.87164572975347
$ws.Range("M8").Value = 16.81322503809883
$ws.Range("B9").Value = 16.79929562856825
$ws.Range("C9").Value = 6.628574613248582
$ws.Range("D9").Value = 11.64874693562948
$ws.Range("F9").Value = 59.2801959350025
$ws.Range("G9").Value = 3.762253990430427
$ws.Range("K9").Value = 12.82754892435429
$ws.Range("L9").Value = 10.86955774124412
$ws.Range("M9").Value = 16.77268928458779
$ws.Range("B10").Value = 16.90971300251965
$ws.Range("C10").Value = 6.755054623188234
$ws.Range("D10").Value = 11.88943589168693
$ws.Range("F10").Value = 61.22580887005954
$ws.Range("G10").Value = 3.756647817077018
$ws.Range("K10").Value = 12.85076743675962
$ws.Range("L10").Value = 10.87533761663739
$ws.Range("M10").Value = 16.76064529309658
$ws.Range("B11").Value = 16.9662154915052
$ws.Range("C11").Value = 6.810280818901685
$ws.Range("D11").Value = 11.99632700244635
$ws.Range("F11").Value = 62.09272834823341
$ws.Range("G11").Value = 3.754211047246488
$ws.Range("K11").Value = 12.86680009322671
$ws.Range("L11").Value = 10.87954877993519
$ws.Range("M11").Value = 16.75901918638848
$ws.Range("B12").Value = 16.9884992228175
$ws.Range("C12").Value = 6.830862772309226
$ws.Range("D12").Value = 12.03641746157614
$ws.Range("F12").Value = 62.41817907781721
$ws.Range("G12").Value = 3.75330451106632
$ws.Range("K12").Value = 12.87365564472956
$ws.Range("L12").Value = 10.88137014950083
$ws.Range("M12").Value = 16.75895701100678
$ws.Range("B13").Value = 16.98366080309127
$ws.Range("C13").Value = 6.826444782517586
$ws.Range("D13").Value = 12.02780067375486
$ws.Range("F13").Value = 62.34821683356174
$ws.Range("G13").Value = 3.753499030347291
$ws.Range("K13").Value = 12.87214434503819
$ws.Range("L13").Value = 10.88096781624978
$ws.Range("M13").Value = 16.75894579030614
$ws.Range("B14").Value = 16.96803109112886
$ws.Range("C14").Value = 6.8119807354663
$ws.Range("D14").Value = 11.9996330955631
$ws.Range("F14").Value = 62.11956137774015
$ws.Range("G14").Value = 3.754136141591546
$ws.Range("K14").Value = 12.86734840183598
$ws.Range("L14").Value = 10.87969408666645
$ws.Range("M14").Value = 16.75900297983009
$ws.Range("B15").Value = 16.95857255259336
$ws.Range("C15").Value = 6.803078018293172
$ws.Range("D15").Value = 11.98232887696849
$ws.Range("F15").Value = 61.97912800806993
$ws.Range("G15").Value = 3.754528499230463
$ws.Range("K15").Value = 12.8645128025626
$ws.Range("L15").Value = 10.87894338778672
$ws.Range("M15").Value = 16.7591100867022
$ws.Range("B16").Value = 16.9061454545373
$ws.Range("C16").Value = 6.751399075372487
$ws.Range("D16").Value = 11.88239686845916
$ws.Range("F16").Value = 61.16876904773292
$ws.Range("G16").Value = 3.756809340652349
$ws.Range("K16").Value = 12.84982957381756
$ws.Range("L16").Value = 10.87509417446063
$ws.Range("M16").Value = 16.76082906094688
$ws.Range("B17").Value = 16.87557998637741
$ws.Range("C17").Value = 6.719104855533608
$ws.Range("D17").Value = 11.82041656443039
$ws.Range("F17").Value = 60.66682845507967
$ws.Range("G17").Value = 3.758237558902079
$ws.Range("K17").Value = 12.84222173267571
$ws.Range("L17").Value = 10.87313744234205
$ws.Range("M17").Value = 16.76287017445348
$ws.Range("B18").Value = 16.85859097660417
$ws.Range("C18").Value = 6.700312553195547
$ws.Range("D18").Value = 11.78452320777141
$ws.Range("F18").Value = 60.37642963161198
$ws.Range("G18").Value = 3.759069722383229
$ws.Range("K18").Value = 12.8383610212295
$ws.Range("L18").Value = 10.87216092814374
$ws.Range("M18").Value = 16.76440685438116
$ws.Range("B19").Value = 16.85294078591832
$ws.Range("C19").Value = 6.693912407508495
$ws.Range("D19").Value = 11.7723288460545
$ws.Range("F19").Value = 60.27782133491043
$ws.Range("G19").Value = 3.759353317924245
$ws.Range("K19").Value = 12.83714237346446
$ws.Range("L19").Value = 10.87185590081104
$ws.Range("M19").Value = 16.76498944734565
$ws.Range("B20").Value = 16.87877262128266
$ws.Range("C20").Value = 6.722565134358388
$ws.Range("D20").Value = 11.82703981504945
$ws.Range("F20").Value = 60.72043806380545
$ws.Range("G20").Value = 3.758084417021323
$ws.Range("K20").Value = 12.84297830076292
$ws.Range("L20").Value = 10.87333032973843
$ws.Range("M20").Value = 16.76261536009286
$ws.Range("B21").Value = 16.97259795575443
$ws.Range("C21").Value = 6.816238150397963
$ws.Range("D21").Value = 12.00791719782212
$ws.Range("F21").Value = 62.18680161833855
$ws.Range("G21").Value = 3.753948567333519
$ws.Range("K21").Value = 12.86873582433971
$ws.Range("L21").Value = 10.88006206625532
$ws.Range("M21").Value = 16.75897116236364
$ws.Range("B22").Value = 17.03908308900852
$ws.Range("C22").Value = 6.875530102047624
$ws.Range("D22").Value = 12.1238702472116
$ws.Range("F22").Value = 63.12854937412538
$ws.Range("G22").Value = 3.751340018640858
$ws.Range("K22").Value = 12.89013963751932
$ws.Range("L22").Value = 10.88578255143445
$ws.Range("M22").Value = 16.75981583463902
$ws.Range("B23").Value = 17.00313140214564
$ws.Range("C23").Value = 6.84406094708226
$ws.Range("D23").Value = 12.06219497709018
$ws.Range("F23").Value = 62.62751032023648
$ws.Range("G23").Value = 3.752723641410014
$ws.Range("K23").Value = 12.87829895337368
$ws.Range("L23").Value = 10.88260884241691
$ws.Range("M23").Value = 16.75907002460926
$ws.Range("B24").Value = 16.87732741244985
$ws.Range("C24").Value = 6.721001446167428
$ws.Range("D24").Value = 11.82404625254036
$ws.Range("F24").Value = 60.69620683170175
$ws.Range("G24").Value = 3.758153618000577
$ws.Range("K24").Value = 12.84263465777972
$ws.Range("L24").Value = 10.87324266293395
$ws.Range("M24").Value = 16.76272943025028
$ws.Range("B25").Value = 16.76428138737838
$ws.Range("C25").Value = 6.579940854048247
$ws.Range("D25").Value = 11.55789496418045
$ws.Range("F25").Value = 58.54916744018045
$ws.Range("G25").Value = 3.764419879845739
$ws.Range("K25").Value = 12.82385706584375
$ws.Range("L25").Value = 10.86883627466116
$ws.Range("M25").Value = 16.78054163316884
